$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows, shifting existing data down, to match the
# target layout (row numbers below are the CURRENT sheet position at the
# moment of each insert, since earlier inserts shift later rows down):
#   - new blank row above current row 4  (old row4  -> row5,  ...)
#   - new blank row above current row 6  (old row5, now at row6, -> row7, ...)
#   - new blank row above current row 15 (old row13, now at row15, -> row16, ...)
$ws.Rows.Item(4).EntireRow.Insert()
$ws.Rows.Item(6).EntireRow.Insert()
$ws.Rows.Item(15).EntireRow.Insert()

# Update the active selection to match the new layout (row 15 is now the
# newly inserted blank row just above the last two data rows).
$ws.Range("A15:XFD15").Select()
